{"js": "// Replace the date line and the 25 \"a\u00f7b=\" division prompts with their\n// updated values, per the authoring diff. Each (old, new) pair below is\n// unique within the document, so a simple case-sensitive search + full\n// replace of the matched range is safe and preserves the existing run\n// formatting (font/size) because insertText(\"Replace\") rewrites the text\n// of the already-formatted run/range in place.\nconst replacements = [\n  [\"2023-11-04 Saturday\", \"2023-11-05 Sunday\"],\n  [\"84\u00f73=\", \"77\u00f77=\"],\n  [\"28\u00f78=\", \"34\u00f79=\"],\n  [\"33\u00f76=\", \"38\u00f78=\"],\n  [\"94\u00f73=\", \"74\u00f74=\"],\n  [\"51\u00f74=\", \"74\u00f79=\"],\n  [\"25\u00f72=\", \"72\u00f79=\"],\n  [\"16\u00f78=\", \"61\u00f76=\"],\n  [\"54\u00f77=\", \"56\u00f74=\"],\n  [\"83\u00f75=\", \"71\u00f75=\"],\n  [\"24\u00f74=\", \"35\u00f74=\"],\n  [\"78\u00f73=\", \"52\u00f73=\"],\n  [\"51\u00f75=\", \"60\u00f75=\"],\n  [\"82\u00f72=\", \"64\u00f75=\"],\n  [\"84\u00f72=\", \"44\u00f75=\"],\n  [\"94\u00f78=\", \"89\u00f72=\"],\n  [\"21\u00f76=\", \"14\u00f74=\"],\n  [\"52\u00f77=\", \"97\u00f73=\"],\n  [\"37\u00f72=\", \"17\u00f78=\"],\n  [\"30\u00f72=\", \"33\u00f77=\"],\n  [\"45\u00f72=\", \"88\u00f78=\"],\n  [\"78\u00f74=\", \"59\u00f78=\"],\n  [\"35\u00f76=\", \"48\u00f76=\"],\n  [\"87\u00f78=\", \"96\u00f74=\"],\n  [\"13\u00f72=\", \"20\u00f77=\"],\n  [\"70\u00f76=\", \"65\u00f79=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 \"a/b=\" division prompts with their\n# updated values, per the authoring diff. Each Old value is unique in the\n# document, so Find/Replace (wdReplaceAll, but really only ever one hit)\n# safely rewrites the text of the already-formatted run in place without\n# touching font/size formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2023-11-04 Saturday\"; New = \"2023-11-05 Sunday\" },\n    @{ Old = \"84\u00f73=\"; New = \"77\u00f77=\" },\n    @{ Old = \"28\u00f78=\"; New = \"34\u00f79=\" },\n    @{ Old = \"33\u00f76=\"; New = \"38\u00f78=\" },\n    @{ Old = \"94\u00f73=\"; New = \"74\u00f74=\" },\n    @{ Old = \"51\u00f74=\"; New = \"74\u00f79=\" },\n    @{ Old = \"25\u00f72=\"; New = \"72\u00f79=\" },\n    @{ Old = \"16\u00f78=\"; New = \"61\u00f76=\" },\n    @{ Old = \"54\u00f77=\"; New = \"56\u00f74=\" },\n    @{ Old = \"83\u00f75=\"; New = \"71\u00f75=\" },\n    @{ Old = \"24\u00f74=\"; New = \"35\u00f74=\" },\n    @{ Old = \"78\u00f73=\"; New = \"52\u00f73=\" },\n    @{ Old = \"51\u00f75=\"; New = \"60\u00f75=\" },\n    @{ Old = \"82\u00f72=\"; New = \"64\u00f75=\" },\n    @{ Old = \"84\u00f72=\"; New = \"44\u00f75=\" },\n    @{ Old = \"94\u00f78=\"; New = \"89\u00f72=\" },\n    @{ Old = \"21\u00f76=\"; New = \"14\u00f74=\" },\n    @{ Old = \"52\u00f77=\"; New = \"97\u00f73=\" },\n    @{ Old = \"37\u00f72=\"; New = \"17\u00f78=\" },\n    @{ Old = \"30\u00f72=\"; New = \"33\u00f77=\" },\n    @{ Old = \"45\u00f72=\"; New = \"88\u00f78=\" },\n    @{ Old = \"78\u00f74=\"; New = \"59\u00f78=\" },\n    @{ Old = \"35\u00f76=\"; New = \"48\u00f76=\" },\n    @{ Old = \"87\u00f78=\"; New = \"96\u00f74=\" },\n    @{ Old = \"13\u00f72=\"; New = \"20\u00f77=\" },\n    @{ Old = \"70\u00f76=\"; New = \"65\u00f79=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"No match found for '$($pair.Old)'\"\n    }\n}\n"}
